$d = $word.ActiveDocument

# --- Bookmarks around the heading paragraphs --------------------------------
function Add-HeadingBookmark($headingText, $bookmarkName) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs($i)
        if ($p.Range.Text.TrimEnd("`r`a`v") -eq $headingText) {
            $r = $p.Range
            [void]$r.MoveEnd(1, -1)
            $d.Bookmarks.Add($bookmarkName, $r)
            break
        }
    }
}

Add-HeadingBookmark "Mini OpenD6 Legend Rules" "mini-opend6-legend-rules"
Add-HeadingBookmark "Dice Basics" "dice-basics"
Add-HeadingBookmark "Combat" "combat"
Add-HeadingBookmark "Optional Combat Focus Rules" "optional-combat-focus-rules"

# --- Text fix: Target Number -> Target Success conversion rule -------------
# Replace just the final run's text ("divide the TN by 7, rounding up or
# down as you see fit." -> "divide the TN by 6, rounding to the nearest
# number.") while keeping it as its own run (distinct from the preceding
# " " run). A plain Range.Text assignment / Find-Replace would coalesce the
# touched run back into its same-formatted neighbour, so instead delete the
# whole old run (so it disappears entirely) and insert the new text right
# after the gap it leaves -- this creates a brand-new run rather than
# rewriting text inside the merged run set.
$oldText = "divide the TN by 7, rounding up or down as you see fit."
$newText = "divide the TN by 6, rounding to the nearest number."

$full = $d.Content.Text
$idx = $full.IndexOf($oldText)
if ($idx -ge 0) {
    $r = $d.Range($idx, $idx + $oldText.Length)
    $r.Delete()
    $ins = $d.Range($idx, $idx)
    $ins.InsertAfter($newText)
}
